# ----------------------------------------------------------------------
# Edit script: rename "Tabelle2" sheet, refresh its survey-result cells,
# resize rows and restore the view state (zoom/selection) on both sheets.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename the second sheet from the placeholder "Tabelle2" to its real name
$ws2.Name = "Iteration 2 ohne Codefehler"

# 2) Replace the review-output cells on the "without errors" sheet with the
#    real LLM feedback text for the bug-free code sample (rows 2 and 3).
$newText0 = @"
Dein Code ist funktional korrekt und übersichtlich, die Schleifenbedingung „i < nums.Length“ passt.Du könntest stattdessen einmal LINQ einsetzen, z. B.
return nums.Select(x => x + 1).ToArray();
um die Logik als Ein-Zeiler auszudrücken.Wenn du Speicher sparen willst und das Original-Array nicht mehr benötigst, lässt sich auch direkt in-place hochzählen (ohne neues Array).
"@
$ws2.Range("D2").Value = $newText0

$newText1 = @"
Dein Algorithmus ist bereits korrekt und durchläuft das Array in linearer Zeit (O(n)).
• Du könntest mit LINQ etwas kompakter schreiben, etwa nums.Select(n => n + 1).ToArray().
• Wenn keine Kopie nötig ist, ließe sich das Original-Array in-place anpassen und so Speicher sparen.
• Eine expression-bodied Methode (=>) steigert zusätzlich die Lesbarkeit.
"@
$ws2.Range("E2").Value = $newText1

$newText2 = @"
Dein Code erfüllt die Aufgabe korrekt und ist gut lesbar.Du könntest Array.ConvertAll(nums, n => n + 1) oder LINQ (nums.Select(n => n + 1).ToArray()) einsetzen, um den Kern in eine Zeile zu packen.Für extrem große Arrays ließe sich auch Parallel.For nutzen, um die Arbeit auf mehrere Threads zu verteilen.In performancekritischen Szenarien bieten sich Span<T> oder SIMD-Vektorisierungen unter .NET Core an.
"@
$ws2.Range("F2").Value = $newText2

$newText3 = @"
Der Ansatz funktioniert einwandfrei und ist klar verständlich.
Als Alternative könntest du LINQ nutzen, etwa mit nums.Select(n => n + 1).ToArray(), um die Schleife prägnanter auszudrücken.
Ebenso wäre Array.ConvertAll(nums, n => n + 1) eine kompakte Möglichkeit.
Falls dir Performance und Speicherfehler wichtig sind, lohnt sich ein Blick auf Span<T> für zero-allocation-Szenarien.
"@
$ws2.Range("G2").Value = $newText3

$newText4 = @"
Dein Ansatz ist korrekt und durchläuft das Array in O(n), aber hier ein paar Denkanstöße zur Verbesserung:Du könntest auf eine separate Ergebnis-Array verzichten und stattdessen das übergebene Array direkt anpassen, wenn keine ursprüngliche Kopie benötigt wird.Mit LINQ ließe sich die Methode sehr knapp als return nums.Select(x => x + 1).ToArray(); schreiben und wäre zugleich gut lesbar.Für High-Performance-Szenarien könnte ein Span<int> oder Memory<int> helfen, um Zwischenspeicher zu vermeiden.
"@
$ws2.Range("H2").Value = $newText4

$newText5 = @"
Dein Ansatz ist korrekt und durchläuft das Array in O(n), um jedes Element um 1 zu erhöhen.Mit LINQ ließe sich die Schleife kompakter ausdrücken, etwa über nums.Select(n => n + 1).ToArray().Bei sehr großen Arrays könnte man alternativ mit Span<T> oder einer In-Place-Operation experimentieren, um zusätzlichen Speicher zu sparen.Auch Array.ConvertAll(nums, n => n + 1) bietet eine eingebaute, gut lesbare Alternative.
"@
$ws2.Range("I2").Value = $newText5

$newText6 = @"
Dein Code erfüllt die Aufgabe und ist schon recht klar strukturiert.
Überlege, ob du die Zugriffe auf nums.Length vor der Schleife in einer lokalen Variable cachest – ein winziger Performance-Gewinn bei sehr großen Arrays.
Alternativ könntest du mit LINQ etwa nums.Select(x => x + 1).ToArray() einsetzen, um den Code kürzer und deklarativer zu gestalten.
Falls dir ein In-Place-Update genügt, könntest du auch direkt über nums iterieren und die Originalwerte überschreiben.
"@
$ws2.Range("J2").Value = $newText6

$newText7 = @"
Dein Code ist bereits funktional korrekt und läuft in O(n) mit minimalem Speicher-Overhead.
Als Alternative könntest du statt der manuellen Schleife auch LINQ verwenden, z. B.return nums.Select(x => x + 1).ToArray();
Das macht den Code kompakter, wirkt aber intern ähnlich performant.
"@
$ws2.Range("K2").Value = $newText7

$newText8 = @"
Dein Ansatz läuft bereits in O(n) und liefert korrekt für jedes Element eins drauf.Wenn eine In-Place-Änderung ausreicht, könntest du auf das neue Array verzichten und direkt in nums inkrementieren.Für dreizeilige Kürze und höhere Lesbarkeit eignet sich LINQ: nums.Select(n => n + 1).ToArray().Eine expression-bodied Method (public int[] IncrementArray(int[] nums) => nums.Select(n => n+1).ToArray();) macht den Code noch kompakter.
"@
$ws2.Range("L2").Value = $newText8

$newText9 = @"
Der Code ist korrekt lauffähig und sehr übersichtlich.
Um den Zugriff auf nums.Length in jeder Schleifeniteration zu vermeiden, könntest du den Wert einmal in einer lokalen Variable speichern.
Alternativ lässt sich die Aufgabe noch kompakter mit Array.ConvertAll(nums, x => x + 1) oder LINQ (nums.Select(x => x + 1).ToArray()) lösen.
Falls das Original-Array nicht erhalten bleiben muss, käme auch eine In-Place-Änderung mittels Span<int> oder direkter Schleifenmodifikation infrage.
"@
$ws2.Range("M2").Value = $newText9

$newText10 = @"
Hier ein paar kurze Hinweise, wie du deinen Code weiter verbessern kannst:Mögliche Fehler/RisikenNull-Referenz: Wenn nums unerwartet null ist, gibt es eine Exception.Überlauf: Würdest du an einem Element int.MaxValue + 1 rechnen, käme es zum Wrap-around.Nächster sinnvoller SchrittÜberlege, ob du eine Eingangs-Validierung (z. B. auf null) brauchst oder in-place arbeitest, um Speicher zu sparen.Alternativ kannst du mit LINQ in einer Zeile transformieren:
„nums.Select(x => x + 1).ToArray()“.Stil-/Verständlichkeits-VerbesserungenImmer geschweifte Klammern um Schleifen/Ifs setzen, auch bei nur einer Zeile.Einheitliche Einrückung und bei lokalen Variablen evtl. var verwenden.Methodensignatur als Expression‐Bodied Member:public int[] IncrementArray(int[] nums) => nums.Select(x => x + 1).ToArray();
"@
$ws2.Range("G3").Value = $newText10

$newText11 = @"
Hier ein kurzer Hinweis zur Verbesserung – ohne die komplette Lösung vorwegzunehmen:Mögliche Fehler/Risiken
– Null‐Eingabe: Was passiert, wenn nums null ist?
– Integer‐Overflow: nums[i] könnte int.MaxValue sein und beim „+1“ überlaufen.Nächster sinnvoller Schritt
– Füge eine Eingabeprüfung ein (z. B. auf null) und entscheide, wie du damit umgehen willst.
– Überlege, ob du bei Bedarf in einem checked‐Block arbeitest, um Überläufe zu erkennen.Stil-/Verständlichkeits­verbesserungen
– Klammern um den forBlock komplettieren, auch bei einzeiligen Schleifen, für bessere Lesbarkeit.
– Erwäge den Einsatz von LINQ (nums.Select(x => x + 1).ToArray()) oder eines expression-bodied members für noch kompakteren Code.
– Nutze im Methodenkörper ruhig var, wenn der Typ klar ist (z. B. var result = new int[nums.Length];).
"@
$ws2.Range("H3").Value = $newText11

$newText12 = @"
Hier ein paar kurze Hinweise, ohne gleich die Komplett-Lösung vorwegzunehmen:Mögliche Fehler/RisikenWas passiert, wenn jemand nums als null übergibt?Bei sehr großen Werten (z.B. int.MaxValue) kann das Hochzählen einen Überlauf auslösen.Nächster sinn­voller VerbesserungsschrittFüge eine Argument­validierung hinzu (z.B. if (nums == null) throw new ArgumentNullException(...)).Überlege, ob du im überlauf-sensitiven Szenario ein checkedKontext oder zumindest einen Pre-Check auf int.MaxValue brauchst.Kleine Stil- und Verständlichkeits­verbesserungenNutze durchgängig geschweifte Klammern für Schleifen/Blöcke, auch wenn nur eine Zeile folgt.Ziehe var für lokale Variablen in Betracht, wenn der Typ klar ist.Für eine kompakte Variante könntest du später einmal mit LINQ experimentieren (z.B. nums.Select(n => n + 1).ToArray()), beachte aber den Overhead.Beispiel für Null-Check und Stil:public int[] IncrementArray(int[] nums)
{ if (nums == null) throw new ArgumentNullException(nameof(nums)); var result = new int[nums.Length]; for (var i = 0; i < nums.Length; i++) { checked { result[i] = nums[i] + 1; } } return result;
}
"@
$ws2.Range("I3").Value = $newText12

$newText13 = @"
Hier ein paar kurze Hinweise zur Verbesserung deines Codes, ohne gleich die komplette Lösung vorwegzunehmen:Mögliche Fehler/Risiken
– Null-Eingabe: Wenn nums null ist, wirft dein Code eine unhandliche NullReferenceException. Ein früher Guard-Check hilft hier.
– Integer-Overflow: Wenn ein Wert in nums bei int.MaxValue liegt, führt +1 zu einem Überlauf (sofern nicht im checked-Kontext).Nächster sinnvoller Schritt
– Füge einen Guard-Check hinzu:
csharp if (nums == null) throw new ArgumentNullException(nameof(nums)); 
– Entscheide, ob du Überlauf prüfen möchtest (z. B. checkedBlock oder manuelle Prüfung).Stil-/Verständlichkeits-Tipps
– Verwende stets geschweifte Klammern { … } auch bei einzeiligen Schleifen, um Lesbarkeit und Wartbarkeit zu erhöhen.
– Überlege, ob du für so einfache Transformationen nicht LINQ oder Array.ConvertAll nutzen möchtest, um den Code kompakter zu gestalten.Damit hast du eine solide Basis, um die Methode robust und gut lesbar zu machen.
"@
$ws2.Range("J3").Value = $newText13

$newText14 = @"
Hier ein paar Hinweise zur schrittweisen Verbesserung deines Codes:Mögliche Fehler oder RisikenWenn nums null ist, wirft deine Methode eine NullReferenceException. Auch wenn du davon ausgehst, dass immer ein Array reinkommt, lohnt sich ein kurzer Null-Check oder eine klare Fehlermeldung (ArgumentNullException).Nächster sinnvoller VerbesserungsschrittNutze LINQ für mehr Lesbarkeit und Kürze. Statt der klassischen Schleife könntest du schreiben:
return nums.Select(n => n + 1).ToArray();
So beschränkst du dich auf die Geschäftslogik („jeden Wert um 1 erhöhen“) und überlässt die Iteration dem Framework.Stil- und VerständlichkeitsverbesserungenBaue bei Schleifen immer geschweifte Klammern ein, auch wenn sie nur eine Zeile umfassen – das hilft, Fehler beim Hinzufügen von Code zu vermeiden.Erwäge aussagekräftigere Namen (z. B. original statt nums, wenn es die Lesbarkeit verbessert).Mach die Methode statisch, falls sie keinen Objektzustand nutzt:
public static int[] IncrementArray(int[] nums) { … }
"@
$ws2.Range("K3").Value = $newText14

$newText15 = @"
Hier ein paar Punkte, die du noch verbessern kannst, ohne die Lösung vorwegzunehmen:Mögliche Fehler/Risiken
– Was passiert, wenn der Aufrufer null übergibt? In der aktuellen Fassung würdest du eine NullReferenceException bekommen.
– Bei sehr großen Werten könnte das Inkrement theoretisch zum Overflow führen (in deinem Anwendungsbereich aktuell unwahrscheinlich, aber im Idealfall im Blick behalten).Nächster sinnvoller Schritt
– Füge eine Argument-Prüfung hinzu, z. B. if (nums == null) throw new ArgumentNullException(nameof(nums));.
– Überlege, ob du statt der Schleife ein sprachfeatures wie LINQ (nums.Select(...)) oder Array.ConvertAll einsetzen möchtest, um den Code kürzer und deklarativer zu machen.Stil-/Verständlichkeitsverbesserungen
– Auch bei einzelnen Zeilen in einer Schleife sollte man aus Konsistenzgründen immer geschweifte Klammern benutzen.
– Du könntest die Methode als Expression-Bodied Member schreiben, um sie kompakter zu gestalten.Diese Schritte erhöhen Robustheit, Lesbarkeit und nutzen moderne C#-Features.
"@
$ws2.Range("L3").Value = $newText15

$newText16 = @"
Hier ein paar gezielte Hinweise, ohne die komplette Lösung vorwegzunehmen:Mögliche Fehler/RisikenWas passiert, wenn nums null ist? Eine zusätzliche Prüfung (Argument-Validation) verhindert eine NullReferenceException.Könnten Überläufe bei extrem großen Werten auftreten? Im vorliegenden Wertebereich (±10⁹) ist int zwar sicher, aber ein Gedanke für allgemeineren Code.Nächster sinn­voller VerbesserungsschrittÜberlege, ob du statt der manuellen Schleife LINQ einsetzen möchtest. Das macht den Code kurz und ausdrucksstark, ist aber nicht in jedem Szenario performanter.Beispiel als Hinweis (kein vollständiges Rewrite!):
return nums.Select(n => n + 1).ToArray();
Stil- und Verständlichkeits­verbesserungenZieh in Betracht, den Methodenkörper als Expression-Bodied Member zu schreiben, wenn wirklich nur eine Zeile übrig bleibt.Achte auf einheitliche Einrückungen und Klammer­setzung – das erhöht die Lesbarkeit, besonders bei späteren Anpassungen.
"@
$ws2.Range("M3").Value = $newText16

# 3) Row heights settle to new (smaller) values once the long placeholder
#    duplicate text is replaced by the real, differently sized content.
$ws1.Rows.Item(2).RowHeight = 255
$ws1.Rows.Item(4).RowHeight = 195
$ws2.Rows.Item(2).RowHeight = 300
$ws2.Rows.Item(4).RowHeight = 285

# 4) Restore each sheet's zoom level + last selected cell.
#    "Iteration 2 mit Codefehlern" (sheet 1) view state:
$ws1.Activate()
$excel.ActiveWindow.Zoom = 55
$ws1.Range("K3").Select()

#    "Iteration 2 ohne Codefehler" (sheet 2) view state - also the tab that
#    stays active/selected when the workbook is reopened.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 55
$ws2.Range("N2").Select()

Write-Output "edit complete"
